$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.136.99"
$ws.Range("E2").Value = "  +0.96%  "

$ws.Range("D3").Value = "'1.891.87"
$ws.Range("E3").Value = "  +0.26%  "

$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "'0.7423"
$ws.Range("E5").Value = "  -0.50%  "

$ws.Range("D6").Value = "'242.71"
$ws.Range("E6").Value = "  +0.08%  "

$ws.Range("D7").Value = "'1.002"
$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("D8").Value = "'0.3169"
$ws.Range("E8").Value = "  +1.83%  "

$ws.Range("D9").Value = "'0.07214"
$ws.Range("E9").Value = "  +1.41%  "

$ws.Range("D10").Value = "'24.96"
$ws.Range("E10").Value = "  -1.26%  "

$ws.Range("D11").Value = "'0.08354"
$ws.Range("E11").Value = "  -1.64%  "

$ws.Range("D12").Value = "'1.954.25"
$ws.Range("E12").Value = "  +2.74%  "

$ws.Range("D13").Value = "'0.7611"
$ws.Range("E13").Value = "  +0.20%  "

$ws.Range("D14").Value = "'5.436"
$ws.Range("E14").Value = "  +1.56%  "

$ws.Range("D15").Value = "'92.88"
$ws.Range("E15").Value = "  -0.48%  "

$ws.Range("D16").Value = "'6.182"
$ws.Range("E16").Value = "  +0.51%  "

$ws.Range("D17").Value = "'30.199.71"
$ws.Range("E17").Value = "  +1.02%  "

$ws.Range("D18").Value = "'250.27"
$ws.Range("E18").Value = "  +2.93%  "

$ws.Range("D19").Value = "'13.65"
$ws.Range("E19").Value = "  -0.26%  "

$ws.Range("D20").Value = "'0.000007864"
$ws.Range("E20").Value = "  +0.96%  "

$ws.Range("D21").Value = "'2.168.26"
$ws.Range("E21").Value = "  +0.65%  "

$ws.Range("E22").Value = "  +0.05%  "

$ws.Range("D23").Value = "'8.012"
$ws.Range("E23").Value = "  +0.38%  "

$ws.Range("D24").Value = "'1.001"
$ws.Range("E24").Value = "  -0.04%  "

$ws.Range("D25").Value = "'0.1583"
$ws.Range("E25").Value = "  -0.14%  "

$ws.Range("D26").Value = "'9.296"
$ws.Range("E26").Value = "  -0.83%  "

$ws.Range("D27").Value = "'164.06"
$ws.Range("E27").Value = "  +1.13%  "

$ws.Range("D28").Value = "'18.77"
$ws.Range("E28").Value = "  +0.18%  "

$ws.Range("D29").Value = "'2.067"
$ws.Range("E29").Value = "  +2.05%  "

$ws.Range("D30").Value = "'1.479"
$ws.Range("E30").Value = "  -1.60%  "

$ws.Range("D31").Value = "'4.602"
$ws.Range("E31").Value = "  +3.04%  "

$ws.Range("D32").Value = "'1.536"
$ws.Range("E32").Value = "  +0.42%  "

$ws.Range("D33").Value = "'4.218"
$ws.Range("E33").Value = "  +2.86%  "

$ws.Range("D34").Value = "'0.05395"
$ws.Range("E34").Value = "  +0.11%  "

$ws.Range("D35").Value = "'1.251"
$ws.Range("E35").Value = "  +1.05%  "

$ws.Range("D36").Value = "'0.7681"
$ws.Range("E36").Value = "  +3.32%  "

$ws.Range("D37").Value = "'1.002"
$ws.Range("E37").Value = "  +0.05%  "

$ws.Range("D38").Value = "'2.723"
$ws.Range("E38").Value = "  +0.48%  "

$ws.Range("D39").Value = "'0.01971"
$ws.Range("E39").Value = "  +2.03%  "

$ws.Range("D40").Value = "'2.767"
$ws.Range("E40").Value = "  -0.01%  "

$ws.Range("D41").Value = "'0.4558"
$ws.Range("E41").Value = "  +2.40%  "

$ws.Range("D42").Value = "'1.102.13"
$ws.Range("E42").Value = "  +1.33%  "

$ws.Range("D43").Value = "'6.074"
$ws.Range("E43").Value = "  +0.31%  "

$ws.Range("D44").Value = "'72.64"
$ws.Range("E44").Value = "  -0.06%  "

$ws.Range("D45").Value = "'0.8693"
$ws.Range("E45").Value = "  +0.99%  "

$ws.Range("D46").Value = "'104.38"
$ws.Range("E46").Value = "  +1.86%  "

$ws.Range("D47").Value = "'1.002"
$ws.Range("E47").Value = "  +0.13%  "

$ws.Range("D48").Value = "'1.868"
$ws.Range("E48").Value = "  +0.42%  "

$ws.Range("D49").Value = "'7.617"
$ws.Range("E49").Value = "  -0.61%  "

$ws.Range("D50").Value = "'9.620"
$ws.Range("E50").Value = "  -0.59%  "

$ws.Range("D51").Value = "'2.065.36"
$ws.Range("E51").Value = "  +0.55%  "
